$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 3.034777666666667
$ws.Range("H2").Value = 9.104333
$ws.Range("I2").Value = 0.2502264227183869
$ws.Range("J2").Value = 0.2502264227183869
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 10.23061133333333
$ws.Range("N2").Value = 30.691834
$ws.Range("O2").Value = 0.4855635428718841
$ws.Range("P2").Value = 0.4855635428718841
$ws.Range("Q2").Value = 31.04763079074689
$ws.Range("R2").Value = 279.428677116722
$ws.Range("S2").Value = 0.1215008283352977
$ws.Range("T2").Value = 0.1215008283352977

$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 3.034777666666667
$ws.Range("H3").Value = 9.104333
$ws.Range("I3").Value = 0.2502264227183869
$ws.Range("J3").Value = 0.2502264227183869
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 8.775186333333332
$ws.Range("N3").Value = 26.325559
$ws.Range("O3").Value = 0.4164864079521221
$ws.Range("P3").Value = 0.4164864079521222
$ws.Range("Q3").Value = 26.63073950523856
$ws.Range("R3").Value = 239.676655547147
$ws.Range("S3").Value = 0.1042159039726902
$ws.Range("T3").Value = 0.1042159039726903

$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 3.034777666666667
$ws.Range("H4").Value = 9.104333
$ws.Range("I4").Value = 0.2502264227183869
$ws.Range("J4").Value = 0.2502264227183869
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 2.034752
$ws.Range("N4").Value = 6.104255999999999
$ws.Range("O4").Value = 0.09657305490303886
$ws.Range("P4").Value = 0.09657305490303887
$ws.Range("Q4").Value = 6.175019926805333
$ws.Range("R4").Value = 55.575179341248
$ws.Range("S4").Value = 0.02416513005937379
$ws.Range("T4").Value = 0.02416513005937379

$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 3.034777666666667
$ws.Range("H5").Value = 9.104333
$ws.Range("I5").Value = 0.2502264227183869
$ws.Range("J5").Value = 0.2502264227183869
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.02901266666666667
$ws.Range("N5").Value = 0.087038
$ws.Range("O5").Value = 0.001376994272954919
$ws.Range("P5").Value = 0.001376994272954919
$ws.Range("Q5").Value = 0.08804699285044446
$ws.Range("R5").Value = 0.7924229356540001
$ws.Range("S5").Value = 0.0003445603510252153
$ws.Range("T5").Value = 0.0003445603510252153

$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 5.864004666666666
$ws.Range("H6").Value = 17.592014
$ws.Range("I6").Value = 0.4835045831069426
$ws.Range("J6").Value = 0.4835045831069426
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 10.23061133333333
$ws.Range("N6").Value = 30.691834
$ws.Range("O6").Value = 0.4855635428718841
$ws.Range("P6").Value = 0.4855635428718841
$ws.Range("Q6").Value = 59.99235260151956
$ws.Range("R6").Value = 539.931173413676
$ws.Range("S6").Value = 0.2347721983682004
$ws.Range("T6").Value = 0.2347721983682004

$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 5.864004666666666
$ws.Range("H7").Value = 17.592014
$ws.Range("I7").Value = 0.4835045831069426
$ws.Range("J7").Value = 0.4835045831069426
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 8.775186333333332
$ws.Range("N7").Value = 26.325559
$ws.Range("O7").Value = 0.4164864079521221
$ws.Range("P7").Value = 0.4164864079521222
$ws.Range("Q7").Value = 51.45773360953621
$ws.Range("R7").Value = 463.119602485826
$ws.Range("S7").Value = 0.2013730870465988
$ws.Range("T7").Value = 0.2013730870465988

$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 5.864004666666666
$ws.Range("H8").Value = 17.592014
$ws.Range("I8").Value = 0.4835045831069426
$ws.Range("J8").Value = 0.4835045831069426
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 2.034752
$ws.Range("N8").Value = 6.104255999999999
$ws.Range("O8").Value = 0.09657305490303886
$ws.Range("P8").Value = 0.09657305490303887
$ws.Range("Q8").Value = 11.93179522350933
$ws.Range("R8").Value = 107.386157011584
$ws.Range("S8").Value = 0.04669351465025768
$ws.Range("T8").Value = 0.04669351465025769

$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 5.864004666666666
$ws.Range("H9").Value = 17.592014
$ws.Range("I9").Value = 0.4835045831069426
$ws.Range("J9").Value = 0.4835045831069426
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.02901266666666667
$ws.Range("N9").Value = 0.087038
$ws.Range("O9").Value = 0.001376994272954919
$ws.Range("P9").Value = 0.001376994272954919
$ws.Range("Q9").Value = 0.1701304127257778
$ws.Range("R9").Value = 1.531173714532
$ws.Range("S9").Value = 0.0006657830418857155
$ws.Range("T9").Value = 0.0006657830418857155

$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 2.564975
$ws.Range("H10").Value = 7.694925
$ws.Range("I10").Value = 0.2114897989601526
$ws.Range("J10").Value = 0.2114897989601526
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 10.23061133333333
$ws.Range("N10").Value = 30.691834
$ws.Range("O10").Value = 0.4855635428718841
$ws.Range("P10").Value = 0.4855635428718841
$ws.Range("Q10").Value = 26.24126230471667
$ws.Range("R10").Value = 236.17136074245
$ws.Range("S10").Value = 0.1026917360643542
$ws.Range("T10").Value = 0.1026917360643542

$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 2.564975
$ws.Range("H11").Value = 7.694925
$ws.Range("I11").Value = 0.2114897989601526
$ws.Range("J11").Value = 0.2114897989601526
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 8.775186333333332
$ws.Range("N11").Value = 26.325559
$ws.Range("O11").Value = 0.4164864079521221
$ws.Range("P11").Value = 0.4164864079521222
$ws.Range("Q11").Value = 22.50813356534166
$ws.Range("R11").Value = 202.573202088075
$ws.Range("S11").Value = 0.08808262668743042
$ws.Range("T11").Value = 0.08808262668743043

$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 2.564975
$ws.Range("H12").Value = 7.694925
$ws.Range("I12").Value = 0.2114897989601526
$ws.Range("J12").Value = 0.2114897989601526
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 2.034752
$ws.Range("N12").Value = 6.104255999999999
$ws.Range("O12").Value = 0.09657305490303886
$ws.Range("P12").Value = 0.09657305490303887
$ws.Range("Q12").Value = 5.219088011199999
$ws.Range("R12").Value = 46.97179210079999
$ws.Range("S12").Value = 0.02042421596641147
$ws.Range("T12").Value = 0.02042421596641147

$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 2.564975
$ws.Range("H13").Value = 7.694925
$ws.Range("I13").Value = 0.2114897989601526
$ws.Range("J13").Value = 0.2114897989601526
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 0.6666666666666666
$ws.Range("M13").Value = 0.02901266666666667
$ws.Range("N13").Value = 0.087038
$ws.Range("O13").Value = 0.001376994272954919
$ws.Range("P13").Value = 0.001376994272954919
$ws.Range("Q13").Value = 0.07441676468333334
$ws.Range("R13").Value = 0.66975088215
$ws.Range("S13").Value = 0.0002912202419565172
$ws.Range("T13").Value = 0.0002912202419565172

$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 0.6643690000000001
$ws.Range("H14").Value = 1.993107
$ws.Range("I14").Value = 0.05477919521451775
$ws.Range("J14").Value = 0.05477919521451775
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 10.23061133333333
$ws.Range("N14").Value = 30.691834
$ws.Range("O14").Value = 0.4855635428718841
$ws.Range("P14").Value = 0.4855635428718841
$ws.Range("Q14").Value = 6.796901020915334
$ws.Range("R14").Value = 61.172109188238
$ws.Range("S14").Value = 0.0265987801040318
$ws.Range("T14").Value = 0.0265987801040318

$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 0.6643690000000001
$ws.Range("H15").Value = 1.993107
$ws.Range("I15").Value = 0.05477919521451775
$ws.Range("J15").Value = 0.05477919521451775
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 8.775186333333332
$ws.Range("N15").Value = 26.325559
$ws.Range("O15").Value = 0.4164864079521221
$ws.Range("P15").Value = 0.4164864079521222
$ws.Range("Q15").Value = 5.829961769090334
$ws.Range("R15").Value = 52.469655921813
$ws.Range("S15").Value = 0.02281479024540257
$ws.Range("T15").Value = 0.02281479024540257

$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 0.6643690000000001
$ws.Range("H16").Value = 1.993107
$ws.Range("I16").Value = 0.05477919521451775
$ws.Range("J16").Value = 0.05477919521451775
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 2.034752
$ws.Range("N16").Value = 6.104255999999999
$ws.Range("O16").Value = 0.09657305490303886
$ws.Range("P16").Value = 0.09657305490303887
$ws.Range("Q16").Value = 1.351826151488
$ws.Range("R16").Value = 12.166435363392
$ws.Range("S16").Value = 0.005290194226995906
$ws.Range("T16").Value = 0.005290194226995906

$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 0.6643690000000001
$ws.Range("H17").Value = 1.993107
$ws.Range("I17").Value = 0.05477919521451775
$ws.Range("J17").Value = 0.05477919521451775
$ws.Range("K17").Value = 2
$ws.Range("L17").Value = 0.6666666666666666
$ws.Range("M17").Value = 0.02901266666666667
$ws.Range("N17").Value = 0.087038
$ws.Range("O17").Value = 0.001376994272954919
$ws.Range("P17").Value = 0.001376994272954919
$ws.Range("Q17").Value = 0.01927511634066667
$ws.Range("R17").Value = 0.173476047066
$ws.Range("S17").Value = 0.000075430638087470422339236853
$ws.Range("T17").Value = 0.000075430638087470422339236853
